$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 now-obsolete quarters (2007-10-01 .. 2009-10-01) from the
# top of the naive-forecast diff table; this shifts the whole data block
# up by 5 rows (and garbage-collects the now-unused shared strings).
$ws.Rows("2:6").Delete() | Out-Null

# Extend the staircase of forecast-error values: each row now reaches
# one column further to the right than before (bug fix for the naive
# component forecaster).
$ws.Cells.Item(35,2).Value = 0.9944388281790735
$ws.Cells.Item(35,3).Value = -0.7367179991959909
$ws.Cells.Item(35,4).Value = 1.689029222266171
$ws.Cells.Item(35,5).Value = 0.4578386464905256
$ws.Cells.Item(35,6).Value = 2.000627817329037
$ws.Cells.Item(35,7).Value = 1.640654201609292
$ws.Cells.Item(35,8).Value = 2.482106051248113
$ws.Cells.Item(35,9).Value = -0.1934093249009629
$ws.Cells.Item(35,10).Value = 1.04033928895951
$ws.Cells.Item(35,11).Value = 2.764463663985141
$ws.Cells.Item(36,2).Value = -0.7613146397174311
$ws.Cells.Item(36,3).Value = 1.738748518289537
$ws.Cells.Item(36,4).Value = 0.3289336411895696
$ws.Cells.Item(36,5).Value = 1.917178009547636
$ws.Cells.Item(36,6).Value = 1.571442225086492
$ws.Cells.Item(36,7).Value = 2.400889081562568
$ws.Cells.Item(36,8).Value = -0.2730707983855708
$ws.Cells.Item(36,9).Value = 0.9621530155802016
$ws.Cells.Item(36,10).Value = 2.68554092831471
$ws.Cells.Item(36,11).Value = 3.052286043054687
$ws.Cells.Item(37,2).Value = 1.450140857813076
$ws.Cells.Item(37,3).Value = 0.2562259180950384
$ws.Cells.Item(37,4).Value = 1.97433473189567
$ws.Cells.Item(37,5).Value = 1.547668835674672
$ws.Cells.Item(37,6).Value = 2.378006560358896
$ws.Cells.Item(37,7).Value = -0.281103352422687
$ws.Cells.Item(37,8).Value = 0.9493339525436422
$ws.Cells.Item(37,9).Value = 2.671442079878515
$ws.Cells.Item(37,10).Value = 3.039476375130951
$ws.Cells.Item(37,11).Value = -2.436065326563203
$ws.Cells.Item(38,2).Value = 0.6406334362307073
$ws.Cells.Item(38,3).Value = 2.072136993074429
$ws.Cells.Item(38,4).Value = 1.365128770604209
$ws.Cells.Item(38,5).Value = 2.339092923011738
$ws.Cells.Item(38,6).Value = -0.3079149205996747
$ws.Cells.Item(38,7).Value = 0.8899316204641817
$ws.Cells.Item(38,8).Value = 2.619680973924259
$ws.Cells.Item(38,9).Value = 2.99190035428146
$ws.Cells.Item(38,10).Value = -2.486476105632129
$ws.Cells.Item(38,11).Value = 0.8740166842005405
$ws.Cells.Item(39,2).Value = 2.279041869813757
$ws.Cells.Item(39,3).Value = 1.364660634097036
$ws.Cells.Item(39,4).Value = 2.233850469905869
$ws.Cells.Item(39,5).Value = -0.3374747896449133
$ws.Cells.Item(39,6).Value = 0.8572197922478608
$ws.Cells.Item(39,7).Value = 2.572668757538748
$ws.Cells.Item(39,8).Value = 2.950126951381641
$ws.Cells.Item(39,9).Value = -2.527041352139664
$ws.Cells.Item(39,10).Value = 0.8320027763378763
$ws.Cells.Item(40,2).Value = 2.035159727951369
$ws.Cells.Item(40,3).Value = 2.420094282127993
$ws.Cells.Item(40,4).Value = -0.5777405448231215
$ws.Cells.Item(40,5).Value = 0.846716050809707
$ws.Cells.Item(40,6).Value = 2.573428919085527
$ws.Cells.Item(40,7).Value = 2.902190984887824
$ws.Cells.Item(40,8).Value = -2.561820069973615
$ws.Cells.Item(40,9).Value = 0.8026396938474522
$ws.Cells.Item(41,2).Value = 2.925451661310402
$ws.Cells.Item(41,3).Value = -0.4190797873159667
$ws.Cells.Item(41,4).Value = 0.6148343974608563
$ws.Cells.Item(41,5).Value = 2.530497674202357
$ws.Cells.Item(41,6).Value = 2.877258180398436
$ws.Cells.Item(41,7).Value = -2.629361321024937
$ws.Cells.Item(41,8).Value = 0.7447744684086398
$ws.Cells.Item(42,2).Value = 0.3866947907724951
$ws.Cells.Item(42,3).Value = 0.8666183061337129
$ws.Cells.Item(42,4).Value = 2.244119497181896
$ws.Cells.Item(42,5).Value = 2.856605817286514
$ws.Cells.Item(42,6).Value = -2.628012363183588
$ws.Cells.Item(42,7).Value = 0.6893704729749799
$ws.Cells.Item(43,2).Value = 0.859102372451746
$ws.Cells.Item(43,3).Value = 2.322448262598703
$ws.Cells.Item(43,4).Value = 2.841688734323599
$ws.Cells.Item(43,5).Value = -2.630009269297866
$ws.Cells.Item(43,6).Value = 0.6999637172787909
$ws.Cells.Item(44,2).Value = 2.518240533879863
$ws.Cells.Item(44,3).Value = 2.915708371652224
$ws.Cells.Item(44,4).Value = -2.72237770473967
$ws.Cells.Item(44,5).Value = 0.6824695392157508
$ws.Cells.Item(45,2).Value = 3.565613746524331
$ws.Cells.Item(45,3).Value = -2.481270880340997
$ws.Cells.Item(45,4).Value = 0.4239743798924789
$ws.Cells.Item(46,2).Value = -1.451990298217711
$ws.Cells.Item(46,3).Value = 0.7853095085029023
$ws.Cells.Item(47,2).Value = 0.5355100695541125
